$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy border/format from row 3 down to row 8 so row 8 reuses the existing bottom-border style (s=6)
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Numeric cells (B column on rows with plain index numbers)
$ws.Range("B4").Value = 17
$ws.Range("B5").Value = 21
$ws.Range("B6").Value = 25
$ws.Range("B9").Value = 5

# Text cells, written in the same column-by-column order as the source edit
$ws.Range("A4").Value = 'SCRIPT/H01P99A/c00a1001.ssb'
$ws.Range("C4").Value = ' I\''m beat! It was a long day!'
$ws.Range("C5").Value = ' We put in a full day today.\nI\''m exhausted.'
$ws.Range("C6").Value = ' We worked hard all day long.\nI\''m tired!'
$ws.Range("C7").Value = ' We should get some sleep now.'
$ws.Range("C8").Value = ' We can be rested up for another\ngood day tomorrow.'
$ws.Range("B7").Value = '31, 35, 39'
$ws.Range("B8").Value = '45, 49, 53'
$ws.Range("D4").Value = ' Как же я устал! Это был\nдолгий день!'
$ws.Range("D5").Value = ' Мы проработали целый день.\nЯ вымотался.'
$ws.Range("D6").Value = ' Мы работали весь день\nнапролёт. Я так устала!'
$ws.Range("D7").Value = ' Нам нужно поспать.'
$ws.Range("D8").Value = ' Отдохнём и завтра начнём\nновый день!'
$ws.Range("E4").Value = ' Ëàë çå ÿ ôòóàì! Üóï áúì\näïìãéê äåîû!'
$ws.Range("E5").Value = ' Íú ðñïñàáïóàìé øåìúê äåîû.\nŸ âúíïóàìòÿ.'
$ws.Range("E6").Value = ' Íú ñàáïóàìé âåòû äåîû\nîàðñïìæó. Ÿ óàë ôòóàìà!'
$ws.Range("E7").Value = ' Îàí îôçîï ðïòðàóû.'
$ws.Range("E8").Value = ' Ïóäïöîæí é èàâóñà îàœîæí\nîïâúê äåîû!'
$ws.Range("C9").Value = '[CN]The next morning…'
$ws.Range("C10").Value = ' Good morning, [hero]!'
$ws.Range("C11").Value = ' Let\''s make it another great day!'
$ws.Range("B10").Value = '68, 72, 76'
$ws.Range("B11").Value = '82, 86, 90'
$ws.Range("D9").Value = '[CN]На следующее утро...'
$ws.Range("D10").Value = ' Доброе утро, [hero]!'
$ws.Range("D11").Value = ' Давай проведём этот день с\nпользой!'
$ws.Range("E9").Value = '[CN]Îà òìåäôýþåå ôóñï...'
$ws.Range("E10").Value = ' Äïáñïå ôóñï, [hero]!'
$ws.Range("E11").Value = ' Äàâàê ðñïâåäæí üóïó äåîû ò\nðïìûèïê!'

# Row heights (rows 4-8, 10-11 use 28.8; row 9 keeps the default)
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(10).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 28.8

# View: zoom + selection
$excel.ActiveWindow.Zoom = 80
$ws.Range("D5").Select() | Out-Null

